# Apply the "Updated as of 3/15/18" edits to the Basic Manual Gantt Chart sheet.
# The task list for the in-progress / not-started meetings was reshuffled:
#  - the "Dr. Engels" kickoff meeting got renamed to include "& Advisors"
#  - "Meeting # 10" -> "Meeting # 11" (its slot shifted, dates unchanged)
#  - a handful of rows were rewritten with new meetings / dates further down
#    the Gantt chart as the schedule progressed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: just a text correction, dates/status stay the same.
$ws.Range("B29").Value = "Meeting with Dr. Engels & Advisors  #1"

# Row 30: rename only ("Meeting # 10" -> "Meeting # 11"); dates unchanged (43170/43170).
$ws.Range("B30").Value = "Meeting # 11"

# Row 31: "Meeting with Advisor # 4" -> "Meeting # 12", date moves 43172 -> 43173.
$ws.Range("B31").Value = "Meeting # 12"
$ws.Range("C31").Value = 43173
$ws.Range("D31").Value = 43173

# Row 32: "Meeting # 11" -> "First Paper Draft (A)", spans 43132-43173.
$ws.Range("B32").Value = "First Paper Draft (A)"
$ws.Range("C32").Value = 43132
$ws.Range("D32").Value = 43173

# Row 33: "Meeting # 12" -> "Online Café Talk (Dan & Tim)", 43179/43179.
$ws.Range("B33").Value = "Online Café Talk (Dan & Tim)"
$ws.Range("C33").Value = 43179
$ws.Range("D33").Value = 43179

# Row 34: "First Paper Draft Review with Advisor" -> "Online Café Talk (Sudip)", 43180/43180.
$ws.Range("B34").Value = "Online Café Talk (Sudip)"
$ws.Range("C34").Value = 43180
$ws.Range("D34").Value = 43180

# Row 35: "First Paper Draft (A)" -> "Meeting to Introduce John ", 43181/43181.
$ws.Range("B35").Value = "Meeting to Introduce John "
$ws.Range("C35").Value = 43181
$ws.Range("D35").Value = 43181

# Row 36: "Meeting with Advisor # 5" stays, but its date slips from 43180-43181 to 43186/43186.
$ws.Range("B36").Value = "Meeting with Advisor # 5"
$ws.Range("C36").Value = 43186
$ws.Range("D36").Value = 43186

# Row 37: "Meeting # 13" (was blank dates) -> "Meeting 13", 43190/43190.
$ws.Range("B37").Value = "Meeting 13"
$ws.Range("C37").Value = 43190
$ws.Range("D37").Value = 43190
# The E37 duration formula was cached as an empty-string result (dates used to
# be blank); re-assigning it forces the engine to drop the stale cached type
# and recompute now that C37/D37 hold real dates.
$ws.Range("E37").Formula = $ws.Range("E37").Formula

# Row 38: "Online Café Talk" (was blank dates) -> "Meeting with Advisor # 6", 43193/43193.
$ws.Range("B38").Value = "Meeting with Advisor # 6"
$ws.Range("C38").Value = 43193
$ws.Range("D38").Value = 43193
$ws.Range("E38").Formula = $ws.Range("E38").Formula

$wb.Save()
